$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.559.49'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.956.85'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''244.01'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").Value = '''58.63'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.377'
$ws.Range("E9").Value = '  +2.63%  '
$ws.Range("D10").Value = '''0.0807'
$ws.Range("E10").Value = '  -4.74%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").Value = '''22.12'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '2.243.68'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '''0.827'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '''13.66'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = '''5.27'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '1.962.46'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '36.415.73'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").Value = '''69.65'
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("D21").Value = '''228.24'
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").Value = '''5.04'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").Value = '''9.24'
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '''160.23'
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").Value = '''19.40'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = '''0.0618'
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("D34").Value = '''4.29'
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("D37").Value = '''3.34'
$ws.Range("E37").Value = '  +9.10%  '
$ws.Range("D38").Value = '''1.77'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = '''5.77'
$ws.Range("E39").Value = '  -9.28%  '
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '''15.98'
$ws.Range("E44").Value = '  -1.13%  '
$ws.Range("D45").Value = '1.361.81'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D47").Value = '''87.73'
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").Value = '''7.11'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").Value = '''2.82'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("D50").Value = '2.134.64'
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '''43.63'
$ws.Range("E51").Value = '  -5.29%  '
